# Add a "Save" column (H) to the s_vals sheet: a header cell matching the
# style already used by the other headers (e.g. G1's bold/centered/bordered
# style), plus a numeric 0 value underneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed H1 with G1's full formatting (copy includes format + value), then
# overwrite the copied value with the new header text so the style/border
# carries over exactly like the other header cells.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data cell beneath the header.
$ws.Range("H2").Value = 0
